$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 600
$ws.Range("J2").Value = 3.5

# Row 3
$ws.Range("G3").Value = 2.38
$ws.Range("J3").Value = 3.4

# Row 4
$ws.Range("F4").Value = 1.65
$ws.Range("H4").Value = 4.4
$ws.Range("J4").Value = 4.1
$ws.Range("K4").Value = 5.2

# Row 5
$ws.Range("G5").Value = 1.47
$ws.Range("H5").Value = 1.09
$ws.Range("J5").Value = 4.2
$ws.Range("Q5").Value = 1.61

# Row 6
$ws.Range("F6").Value = 2.62
$ws.Range("G6").Value = 3.25
$ws.Range("H6").Value = 2.5
$ws.Range("J6").Value = 2.7
$ws.Range("Q6").Value = 1.69

# Row 7
$ws.Range("G7").Value = 2.06
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 3.1

# Row 8
$ws.Range("P8").Value = 2.2
$ws.Range("Q8").Value = 1.48

# Row 9
$ws.Range("F9").Value = 3
$ws.Range("H9").Value = 2.5
$ws.Range("I9").Value = 2.66
$ws.Range("J9").Value = 3.4
$ws.Range("K9").Value = 3.55
$ws.Range("Q9").Value = 2.02
